$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 33 (hunk 0)
$ws.Range("H33").Value = 643.5
$ws.Range("I33").Value = 84.72221999999999
$ws.Range("K33").Value = 84.72221999999999
$ws.Range("M33").Value = 144.27778
# Row 86 (hunk 1)
$ws.Range("H86").Value = 6000
$ws.Range("J86").Value = 6000
$ws.Range("L86").Value = 6000
$ws.Range("N86").Value = -8246
# Row 89 (hunk 2)
$ws.Range("H89").Value = 6000
$ws.Range("J89").Value = 6000
$ws.Range("L89").Value = 30000
$ws.Range("N89").Value = -41232
# Row 106 (hunk 3)
$ws.Range("H106").Value = 32002.75
$ws.Range("I106").Value = 32002.75
$ws.Range("K106").Value = 32002.75
$ws.Range("M106").Value = -31371.75
# Row 113 (hunk 4)
$ws.Range("H113").Value = 2779.1365
$ws.Range("I113").Value = 3049.9167
$ws.Range("J113").Value = 2454.2
$ws.Range("K113").Value = 3049.9167
$ws.Range("L113").Value = 2454.2
$ws.Range("M113").Value = 204.0832999999998
$ws.Range("N113").Value = -8962.200000000001
# Row 137 (hunk 5)
$ws.Range("H137").Value = 2762.4517
$ws.Range("I137").Value = 1490.75
$ws.Range("J137").Value = 4118.933
$ws.Range("K137").Value = 4472.25
$ws.Range("L137").Value = 12356.799
$ws.Range("M137").Value = -1922.25
$ws.Range("N137").Value = -17456.799

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2 (hunk 6)
$ws.Range("H2").Value = 1635
$ws.Range("J2").Value = 1699.75
$ws.Range("L2").Value = 1699.75
$ws.Range("N2").Value = -1925.75
# Row 32 (hunk 7)
$ws.Range("H32").Value = 15312.468
$ws.Range("I32").Value = 6521.4863
$ws.Range("K32").Value = 6521.4863
$ws.Range("M32").Value = -6234.4863
# Row 45 (hunk 8)
$ws.Range("H45").Value = 1771
$ws.Range("I45").Value = 1399
$ws.Range("J45").Value = 1957
$ws.Range("K45").Value = 1399
$ws.Range("L45").Value = 1957
$ws.Range("M45").Value = -1022
$ws.Range("N45").Value = -2711
# Row 110 (hunk 9)
$ws.Range("H110").Value = 10432.5
$ws.Range("I110").Value = 10432.5
$ws.Range("K110").Value = 10432.5
$ws.Range("M110").Value = -8387.5
# Row 116 (hunk 10)
$ws.Range("H116").Value = 1635
$ws.Range("J116").Value = 1699.75
$ws.Range("L116").Value = 1699.75
$ws.Range("N116").Value = -6287.75
# Row 122 (hunk 11)
$ws.Range("H122").Value = 402963.7
$ws.Range("J122").Value = 4530.4443
$ws.Range("L122").Value = 13591.3329
$ws.Range("N122").Value = -18491.3329

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3 (hunk 12)
$ws.Range("H3").Value = 1635
$ws.Range("J3").Value = 1699.75
$ws.Range("L3").Value = 1699.75
$ws.Range("N3").Value = -1927.75
# Row 94 (hunk 13)
$ws.Range("H94").Value = 567.6875
$ws.Range("I94").Value = 478.07144
$ws.Range("K94").Value = 478.07144
$ws.Range("M94").Value = -27.07144

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 7 (hunk 14)
$ws.Range("H7").Value = 65.7
$ws.Range("I7").Value = 81.666664
$ws.Range("J7").Value = 52.636364
$ws.Range("K7").Value = 81.666664
$ws.Range("L7").Value = 52.636364
$ws.Range("M7").Value = 31.333336
$ws.Range("N7").Value = -278.636364
# Row 50 (hunk 15)
$ws.Range("H50").Value = 49399.8
$ws.Range("J50").Value = 54249.75
$ws.Range("L50").Value = 54249.75
$ws.Range("N50").Value = -55499.75
# Row 51 (hunk 16)
$ws.Range("H51").Value = 23499.666
$ws.Range("I51").Value = 500
$ws.Range("J51").Value = 34999.5
$ws.Range("K51").Value = 500
$ws.Range("L51").Value = 34999.5
$ws.Range("M51").Value = 236
$ws.Range("N51").Value = -36471.5
# Row 60 (hunk 17)
$ws.Range("H60").Value = 37999.855
$ws.Range("J60").Value = 48500
$ws.Range("L60").Value = 48500
$ws.Range("N60").Value = -49522
# Row 61 (hunk 18)
$ws.Range("H61").Value = 23499.666
$ws.Range("I61").Value = 500
$ws.Range("J61").Value = 34999.5
$ws.Range("K61").Value = 500
$ws.Range("L61").Value = 34999.5
$ws.Range("M61").Value = -152
$ws.Range("N61").Value = -35695.5
# Row 62 (hunk 19)
$ws.Range("H62").Value = 48679.89
$ws.Range("I62").Value = 4725.4
$ws.Range("J62").Value = 103623
$ws.Range("K62").Value = 4725.4
$ws.Range("L62").Value = 103623
$ws.Range("M62").Value = -4101.4
$ws.Range("N62").Value = -104871
# Row 65 (hunk 20)
$ws.Range("H65").Value = 48679.89
$ws.Range("I65").Value = 4725.4
$ws.Range("J65").Value = 103623
$ws.Range("K65").Value = 23627
$ws.Range("L65").Value = 518115
$ws.Range("M65").Value = -20507
$ws.Range("N65").Value = -524355
# Row 111 (hunk 21)
$ws.Range("H111").Value = 70702
$ws.Range("J111").Value = 70702
$ws.Range("L111").Value = 70702
$ws.Range("N111").Value = -78882
# Row 122 (hunk 22)
$ws.Range("H122").Value = 692.3333
$ws.Range("J122").Value = 999.5
$ws.Range("L122").Value = 2998.5
$ws.Range("N122").Value = -7898.5

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 5 (hunk 23)
$ws.Range("H5").Value = 428.14285
$ws.Range("I5").Value = 399.4
$ws.Range("K5").Value = 1198.2
$ws.Range("M5").Value = -1086.2
# Row 97 (hunk 24)
$ws.Range("H97").Value = 63
$ws.Range("I97").Value = 63
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 189
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = 307
$ws.Range("N97").ClearContents()
# Row 135 (hunk 25)
$ws.Range("H135").Value = 428.14285
$ws.Range("I135").Value = 399.4
$ws.Range("K135").Value = 3594.6
$ws.Range("M135").Value = -1059.6

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 61 (hunk 26)
$ws.Range("H61").Value = 6747.375
$ws.Range("I61").Value = 7746.75
$ws.Range("K61").Value = 7746.75
$ws.Range("M61").Value = -7544.75
# Row 93 (hunk 27)
$ws.Range("H93").Value = 1158.4445
$ws.Range("I93").Value = 963.63635
$ws.Range("K93").Value = 963.63635
$ws.Range("M93").Value = 284.36365
# Row 113 (hunk 28)
$ws.Range("H113").Value = 6747.375
$ws.Range("I113").Value = 7746.75
$ws.Range("K113").Value = 7746.75
$ws.Range("M113").Value = -5576.75
# Row 136 (hunk 29)
$ws.Range("H136").Value = 4379.2
$ws.Range("I136").Value = 4379.2
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 13137.6
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -10587.6
$ws.Range("N136").ClearContents()

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 113 (hunk 30)
$ws.Range("H113").Value = 3134.6667
$ws.Range("I113").Value = 2250.5
$ws.Range("K113").Value = 6751.5
$ws.Range("M113").Value = -4581.5
# Row 122 (hunk 31)
$ws.Range("H122").Value = 1839.5
$ws.Range("J122").Value = 1565
$ws.Range("L122").Value = 4695
$ws.Range("N122").Value = -9595
# Row 132 (hunk 32)
$ws.Range("H132").Value = 1673.3
$ws.Range("I132").Value = 1165.2
$ws.Range("K132").Value = 3495.6
$ws.Range("M132").Value = -965.6000000000004
